# Updates cryptos list figures (price/volume columns) to match the
# latest scrape, and fixes a couple of swapped row orderings plus one
# coin replacement (OKB -> ImmutableX ordering swap, Mantle -> ARBITRUM).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price text (e.g. "209.24", "1.00",
# "0.0000260") that must stay stored as literal text so values such as
# trailing zeros / very small decimals are not mangled by Excel's
# automatic number detection. We force those via a leading apostrophe
# (the same thing Excel does when a user types '123 into a cell).

$ws.Range("D2").Value = "86.920.54"
$ws.Range("E2").Value = "  +5.70%  "
$ws.Range("D3").Value = "3.257.33"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Formula = "'209.24"
$ws.Range("E5").Value = "  -3.67%  "
$ws.Range("D6").Formula = "'624.95"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Formula = "'0.386"
$ws.Range("E7").Value = "  +33.14%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Formula = "'0.651"
$ws.Range("E8").Value = "  +11.37%  "
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").Formula = "'1.00"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "3.249.14"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("D11").Formula = "'0.575"
$ws.Range("E11").Value = "  -4.26%  "
$ws.Range("D12").Formula = "'0.0000260"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Formula = "'0.178"
$ws.Range("E13").Value = "  +7.57%  "
$ws.Range("D14").Formula = "'34.09"
$ws.Range("E14").Value = "  +5.96%  "
$ws.Range("D15").Value = "3.859.38"
$ws.Range("E15").Value = "  +2.67%  "
$ws.Range("D16").Formula = "'5.23"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "86.873.81"
$ws.Range("E17").Value = "  +6.01%  "
$ws.Range("D18").Value = "3.269.84"
$ws.Range("E18").Value = "  +3.23%  "
$ws.Range("D19").Formula = "'14.04"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Formula = "'2.96"
$ws.Range("E20").Value = "  -7.07%  "
$ws.Range("D21").Formula = "'9.00"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Formula = "'430.00"
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("D23").Formula = "'5.34"
$ws.Range("E23").Value = "  +4.08%  "
$ws.Range("D24").Formula = "'7.16"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").Formula = "'12.14"
$ws.Range("E25").Value = "  +8.44%  "
$ws.Range("D26").Formula = "'5.13"
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("D27").Value = "3.447.66"
$ws.Range("E27").Value = "  +2.85%  "
$ws.Range("D28").Formula = "'76.14"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Formula = "'0.0000129"
$ws.Range("E29").Value = "  +5.13%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Formula = "'0.174"
$ws.Range("E31").Value = "  +12.71%  "
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").Formula = "'8.87"
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("D34").Formula = "'546.55"
$ws.Range("E34").Value = "  -5.76%  "
$ws.Range("D35").Formula = "'1.43"
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("D36").Formula = "'1.96"
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("D37").Formula = "'6.80"
$ws.Range("E37").Value = "  +11.44%  "
$ws.Range("D38").Formula = "'0.137"
$ws.Range("E38").Value = "  -10.68%  "
$ws.Range("D39").Formula = "'22.63"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D41").Formula = "'21.73"
$ws.Range("E41").Value = "  +4.33%  "
$ws.Range("D42").Formula = "'0.396"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").Formula = "'2.01"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").Formula = "'2.91"
$ws.Range("E44").Value = "  -5.10%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Formula = "'155.43"
$ws.Range("E46").Value = "  -3.63%  "
$ws.Range("D47").Formula = "'178.14"
$ws.Range("E47").Value = "  -5.06%  "
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D48").Formula = "'1.33"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Formula = "'44.36"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Formula = "'0.627"
$ws.Range("E51").Value = "  -0.62%  "
